$d = $word.ActiveDocument

# 1. Title heading and bold "What we like/don't like" repeated title
$d.Content.Find.Execute("Play Jewel of the Dragon Free - Unique Dragon-Theme Slot Machine", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jewel of the Dragon for Free", 2)

# 2. "Unique grid gameplay" -> "Unique gameplay mechanics with a 5x4 game grid"
$d.Content.Find.Execute("Unique grid gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Unique gameplay mechanics with a 5x4 game grid", 2)

# 3. "Access to free spins with fixed Wild symbols" -> "Free spins with fixed Wild symbols"
$d.Content.Find.Execute("Access to free spins with fixed Wild symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Free spins with fixed Wild symbols", 2)

# 4. "Potential for generous payouts" -> "Potential for generous winnings"
$d.Content.Find.Execute("Potential for generous payouts", $true, $false, $false, $false, $false, $true, 1, $false, "Potential for generous winnings", 2)

# 5. "Graphics are minimalist" -> "Graphics are rather minimalist"
$d.Content.Find.Execute("Graphics are minimalist", $true, $false, $false, $false, $false, $true, 1, $false, "Graphics are rather minimalist", 2)

# 6. "Free spins cannot be re-triggered" -> "Free spins cannot be retriggered during the same mode"
$d.Content.Find.Execute("Free spins cannot be re-triggered", $true, $false, $false, $false, $false, $true, 1, $false, "Free spins cannot be retriggered during the same mode", 2)

# 7. Meta description sentence
$d.Content.Find.Execute("Read our review of Jewel of the Dragon and play free. Stand out from the crowd with unique 5*4 gameplay, fixed Wilds, and potential for big payouts.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Jewel of the Dragon and play for free. Experience unique gameplay mechanics and potential for generous winnings.", 2)
